# Scheduled market-data refresh: updates the cached Universalis price
# columns (currentAveragePrice / currentAveragePriceNQ/HQ) and the
# derived Leve price/profit columns (H:N) for specific leve rows across
# the job sheets. Values below are this run's refreshed snapshot.
$wb = $excel.ActiveWorkbook

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4869.2
$ws.Range("I86").Value = 3116
$ws.Range("J86").Value = 8960
$ws.Range("K86").Value = 3116
$ws.Range("L86").Value = 8960
$ws.Range("M86").Value = -1993
$ws.Range("N86").Value = -11206

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4869.2
$ws.Range("I89").Value = 3116
$ws.Range("J89").Value = 8960
$ws.Range("K89").Value = 15580
$ws.Range("L89").Value = 44800
$ws.Range("M89").Value = -9964
$ws.Range("N89").Value = -56032

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4822.4653
$ws.Range("I32").Value = 4067
$ws.Range("K32").Value = 4067
$ws.Range("M32").Value = -3780

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1464.5476
$ws.Range("I74").Value = 1507.4375
$ws.Range("J74").Value = 1438.1538
$ws.Range("K74").Value = 1507.4375
$ws.Range("L74").Value = 1438.1538
$ws.Range("M74").Value = -633.4375
$ws.Range("N74").Value = -3186.1538

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1464.5476
$ws.Range("I77").Value = 1507.4375
$ws.Range("J77").Value = 1438.1538
$ws.Range("K77").Value = 7537.1875
$ws.Range("L77").Value = 7190.769
$ws.Range("M77").Value = -3169.1875
$ws.Range("N77").Value = -15926.769

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 142859700
$ws.Range("I88").Value = 2968
$ws.Range("J88").Value = 200002380
$ws.Range("K88").Value = 2968
$ws.Range("L88").Value = 200002380
$ws.Range("M88").Value = -2562
$ws.Range("N88").Value = -200003192

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 142859700
$ws.Range("I91").Value = 2968
$ws.Range("J91").Value = 200002380
$ws.Range("K91").Value = 2968
$ws.Range("L91").Value = 200002380
$ws.Range("M91").Value = -1564
$ws.Range("N91").Value = -200005188

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2150.5715
$ws.Range("I86").Value = 1931.3846
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 1931.3846
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -808.3846000000001
$ws.Range("N86").Value = -7246

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2150.5715
$ws.Range("I89").Value = 1931.3846
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 9656.923000000001
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -4040.923000000001
$ws.Range("N89").Value = -36232

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2924.84
$ws.Range("I31").Value = 1352.6086
$ws.Range("J31").Value = 4264.148
$ws.Range("K31").Value = 1352.6086
$ws.Range("L31").Value = 4264.148
$ws.Range("M31").Value = -1057.6086
$ws.Range("N31").Value = -4854.148

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2924.84
$ws.Range("I34").Value = 1352.6086
$ws.Range("J34").Value = 4264.148
$ws.Range("K34").Value = 1352.6086
$ws.Range("L34").Value = 4264.148
$ws.Range("M34").Value = -1150.6086
$ws.Range("N34").Value = -4668.148

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4451.5
$ws.Range("I62").Value = 4449.5
$ws.Range("J62").Value = 4452.5
$ws.Range("K62").Value = 4449.5
$ws.Range("L62").Value = 4452.5
$ws.Range("M62").Value = -3825.5
$ws.Range("N62").Value = -5700.5

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4451.5
$ws.Range("I65").Value = 4449.5
$ws.Range("J65").Value = 4452.5
$ws.Range("K65").Value = 22247.5
$ws.Range("L65").Value = 22262.5
$ws.Range("M65").Value = -19127.5
$ws.Range("N65").Value = -28502.5

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3170
$ws.Range("I122").Value = 2980
$ws.Range("J122").Value = 3233.3333
$ws.Range("K122").Value = 8940
$ws.Range("L122").Value = 9699.999899999999
$ws.Range("M122").Value = -6490
$ws.Range("N122").Value = -14599.9999

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 4545632.5
$ws.Range("I12").Value = 10000108
$ws.Range("J12").Value = 235.66667
$ws.Range("K12").Value = 30000324
$ws.Range("L12").Value = 707.00001
$ws.Range("M12").Value = -30000151
$ws.Range("N12").Value = -1053.00001

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 16666872
$ws.Range("I23").Value = 50000044
$ws.Range("J23").Value = 286
$ws.Range("K23").Value = 150000132
$ws.Range("L23").Value = 858
$ws.Range("M23").Value = -149999897
$ws.Range("N23").Value = -1328

# CUL row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2360
$ws.Range("J32").Value = 3000
$ws.Range("L32").Value = 9000
$ws.Range("N32").Value = -9566

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7877.6665
$ws.Range("J39").Value = 7877.6665
$ws.Range("L39").Value = 23632.9995
$ws.Range("N39").Value = -24220.9995

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3892.56
$ws.Range("J55").Value = 3892.56
$ws.Range("L55").Value = 11677.68
$ws.Range("N55").Value = -12031.68

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3864.818
$ws.Range("J63").Value = 4689.25
$ws.Range("L63").Value = 14067.75
$ws.Range("N63").Value = -15565.75

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 3864.818
$ws.Range("J66").Value = 4689.25
$ws.Range("L66").Value = 42203.25
$ws.Range("N66").Value = -49691.25

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3115.6853
$ws.Range("I68").Value = 5590.091
$ws.Range("J68").Value = 1414.5312
$ws.Range("K68").Value = 16770.273
$ws.Range("L68").Value = 4243.5936
$ws.Range("M68").Value = -15959.273
$ws.Range("N68").Value = -5865.5936

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3115.6853
$ws.Range("I71").Value = 5590.091
$ws.Range("J71").Value = 1414.5312
$ws.Range("K71").Value = 50310.819
$ws.Range("L71").Value = 12730.7808
$ws.Range("M71").Value = -46254.819
$ws.Range("N71").Value = -20842.7808

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1818708.6
$ws.Range("I113").Value = 3333807.5
$ws.Range("J113").Value = 556126.3
$ws.Range("K113").Value = 10001422.5
$ws.Range("L113").Value = 1668378.9
$ws.Range("M113").Value = -9999252.5
$ws.Range("N113").Value = -1672718.9

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15943168
$ws.Range("J131").Value = 16667908
$ws.Range("L131").Value = 50003724
$ws.Range("N131").Value = -50013804

# GSM row 4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 4000
$ws.Range("I4").Value = 4000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 4000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -3888
$ws.Range("N4").ClearContents()

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1738.2222
$ws.Range("I102").Value = 1476.7693
$ws.Range("J102").Value = 2418
$ws.Range("K102").Value = 1476.7693
$ws.Range("L102").Value = 2418
$ws.Range("M102").Value = 145.2307000000001
$ws.Range("N102").Value = -5662

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 45456856
$ws.Range("I68").Value = 2126.9375
$ws.Range("J68").Value = 166669470
$ws.Range("K68").Value = 2126.9375
$ws.Range("L68").Value = 166669470
$ws.Range("M68").Value = -1377.9375
$ws.Range("N68").Value = -166670968

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 45456856
$ws.Range("I71").Value = 2126.9375
$ws.Range("J71").Value = 166669470
$ws.Range("K71").Value = 10634.6875
$ws.Range("L71").Value = 833347350
$ws.Range("M71").Value = -6890.6875
$ws.Range("N71").Value = -833354838

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3986.25
$ws.Range("I82").Value = 1496.8125
$ws.Range("J82").Value = 13944
$ws.Range("K82").Value = 1496.8125
$ws.Range("L82").Value = 13944
$ws.Range("M82").Value = -1135.8125
$ws.Range("N82").Value = -14666

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3986.25
$ws.Range("I85").Value = 1496.8125
$ws.Range("J85").Value = 13944
$ws.Range("K85").Value = 1496.8125
$ws.Range("L85").Value = 13944
$ws.Range("M85").Value = -248.8125
$ws.Range("N85").Value = -16440

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1774.75
$ws.Range("I81").Value = 1449.5
$ws.Range("J81").Value = 2100
$ws.Range("K81").Value = 2899
$ws.Range("L81").Value = 4200
$ws.Range("M81").Value = -1838
$ws.Range("N81").Value = -6322

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1774.75
$ws.Range("I84").Value = 1449.5
$ws.Range("J84").Value = 2100
$ws.Range("K84").Value = 14495
$ws.Range("L84").Value = 21000
$ws.Range("M84").Value = -9191
$ws.Range("N84").Value = -31608
